# Generate Report for Handoff
#
# Refreshes the handoff report: the handoff package id changes from
# 1e5be50e-e4e3-4bc0-9dd9-62b090cc0fa5 to 57ec59dd-e343-473d-b308-789f6f5a3972,
# the xlf content-hash token changes from 0eff5cde03a379ad7b58bedc4abaa0c0d89ff7cc
# to b52f2727d7e0af38f360c666cf7033d0a68b288a, and the associated handoff
# timestamps move forward. Both the cell text and the (separately stored)
# hyperlink display text need updating on every sheet that references the
# old package id.

$wb = $excel.ActiveWorkbook

$oldId = "1e5be50e-e4e3-4bc0-9dd9-62b090cc0fa5"
$newId = "57ec59dd-e343-473d-b308-789f6f5a3972"
$oldHash = "0eff5cde03a379ad7b58bedc4abaa0c0d89ff7cc"
$newHash = "b52f2727d7e0af38f360c666cf7033d0a68b288a"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = "2016-38-14 04:38:38"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
}

# ----- zh-cn sheet -----
$oldZh = "$oldId.$oldHash.zh-cn.xlf"
$newZh = "$newId.$newHash.zh-cn.xlf"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Range("E2").Value = "2016-03-14 04:38:35"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
    if ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newZh
    }
}

# ----- de-de sheet -----
$oldDe = "$oldId.$oldHash.de-de.xlf"
$newDe = "$newId.$newHash.de-de.xlf"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Range("E2").Value = "2016-03-14 04:38:38"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
    if ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newDe
    }
}
